# The "dataset" sheet gets 6 new data rows appended (rows 159-164), each one
# describing the same tensorflow/ranking repository entry found in the very
# last existing row (158) - only the running "id" in column A increases.
#
# Row 158 is copied whole (values + formatting) down into the 6 new rows so
# every column keeps the same look/format as the rest of the table, then the
# id column (A) is corrected to keep counting up (158, 159, ... 163).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 158
$rowsToAdd = 6

$sourceRange = $ws.Range("A" + $lastRow + ":O" + $lastRow)

for ($i = 1; $i -le $rowsToAdd; $i++) {
    $targetRow = $lastRow + $i
    $sourceRange.Copy()
    $ws.Range("A" + $targetRow + ":O" + $targetRow).PasteSpecial(-4104)  # xlPasteAll
    $ws.Cells.Item($targetRow, 1).Value = $lastRow + $i - 1
}

$excel.CutCopyMode = $false
